# Generate Report for Handoff
# Updates the "b.md" row in each sheet to reflect a completed handoff report:
#   - Overview sheet: Status/Date columns for b.md
#   - zh-cn / de-de sheets: Status, Content Duplicate, Latest Handoff File/Datetime,
#     and Error Detail columns for b.md, plus widen the Error Detail column.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-21 00:25:15"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49fbb87f9089cb5b871c6b04188fad9d6d639747/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2214d01d61cdbf6cca0ef791e9cac2788c1cbd87/e2e/b.md."

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-10-21 00:25:04"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-10-21 00:25:15"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
